$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
  "7-6=",
  "90-46=",
  "35+18=",
  "72+9=",
  "17+69=",
  "7+72=",
  "20+77=",
  "39+58=",
  "15+82=",
  "18+33=",
  "56-55=",
  "68-43=",
  "26+20=",
  "30+32=",
  "91-49=",
  "1+35=",
  "74+23=",
  "78-6=",
  "7+83=",
  "77-32=",
  "51+48=",
  "23-21=",
  "6+56=",
  "43-1=",
  "69-0=",
  "21+43=",
  "65-13=",
  "42+3=",
  "51-2=",
  "61+20=",
  "56+11=",
  "99-2=",
  "38-25=",
  "45-2=",
  "30+16=",
  "2+14=",
  "37-8=",
  "7+25=",
  "34+44=",
  "10+75=",
  "0+27=",
  "28+43=",
  "57+26=",
  "86-78=",
  "41-3=",
  "40+50=",
  "43-5=",
  "68-65=",
  "38+61=",
  "60+9=",
  "36+18=",
  "10+88=",
  "50+4=",
  "86-37=",
  "38+23=",
  "24+73=",
  "98-96=",
  "86-13=",
  "78-44=",
  "5+47=",
  "87-0=",
  "49+41=",
  "56-20=",
  "36+18=",
  "57-24=",
  "15-7=",
  "35-16=",
  "82+0=",
  "51-46=",
  "2+35=",
  "68+10=",
  "64-27=",
  "30-3=",
  "79-25=",
  "75-4=",
  "72-11=",
  "51-42=",
  "68-62=",
  "91-51=",
  "14+46=",
  "57+10=",
  "6+77=",
  "92-39=",
  "46+4=",
  "9+72=",
  "7+41=",
  "79-14=",
  "72+9=",
  "59-23=",
  "46+19=",
  "26+28=",
  "4+11=",
  "28+4=",
  "68-21=",
  "84-61=",
  "74+13=",
  "41+33=",
  "95-17=",
  "5+57=",
  "70-26="
)
$idx = 0
for ($r = 1; $r -le 20; $r++) {
  for ($c = 1; $c -le 5; $c++) {
    $cell = $t.Cell($r, $c)
    $cell.Range.Text = $values[$idx]
    $idx = $idx + 1
  }
}
Write-Output "done: $idx cells updated"